# Add "RomaniaFC" and "SlovakiaFC" test-data sheets, cloned from the "UK"
# sheet template (same layout used for every per-market sheet: LCD800 row,
# D-column header cells, printer list, merged title / code-box cells).
#
# New sheets are inserted at the end of the tab strip (after "Spain") and
# populated with the market name (B2) and JIRA/test-case reference (B4) for
# each new country. Slovakia ends up the active/selected sheet, matching the
# commit's final workbookView state.

$wb = $excel.ActiveWorkbook

$uk    = $wb.Worksheets.Item("UK")
$spain = $wb.Worksheets.Item("Spain")

# Clone the UK template sheet twice, placing the copies after Spain, then
# after the first new sheet, so the final tab order is ... Spain, Romania,
# Slovakia.
$uk.Copy($null, $spain)
$romania = $wb.Worksheets.Item("UK (2)")
$romania.Name = "Romania"

$uk.Copy($null, $romania)
$slovakia = $wb.Worksheets.Item("UK (2)")
$slovakia.Name = "Slovakia"

# Market name header cell (B2) for both new sheets.
$romania.Range("B2").Value = "Romania Market"
$slovakia.Range("B2").Value = "Slovakia Market"

# Test-case / ticket reference cell (B4) for both new sheets.
$romania.Range("B4").Value = "NGC-4307/T3536/T3543"
$slovakia.Range("B4").Value = "NGC-4306/T3562/T3575"

# On the Slovakia sheet the reference cell was entered without the usual
# bordered-cell style (unlike every other sheet's B4).
$slovakia.Range("B4").ClearFormats()

# Restore/settle selections: Romania's cursor on B9, Spain's selection back
# to a whole-sheet selection (left over from the copy operation), and finish
# with Slovakia active with B4 selected.
$romania.Range("B9").Select()

$spain.Activate()
$spain.Cells.Select()

$slovakia.Activate()
$slovakia.Range("B4").Select()
